$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Remove the stray helper values that lived in column K (K4, K15,
#    K24, K31, K48, K55, K63, K72). Removing them also shrinks the
#    sheet dimension / row spans back down to column I automatically.
# ---------------------------------------------------------------------
$null = $ws.Range("K4").ClearContents()
$null = $ws.Range("K15").ClearContents()
$null = $ws.Range("K24").ClearContents()
$null = $ws.Range("K31").ClearContents()
$null = $ws.Range("K48").ClearContents()
$null = $ws.Range("K55").ClearContents()
$null = $ws.Range("K63").ClearContents()
$null = $ws.Range("K72").ClearContents()

# ---------------------------------------------------------------------
# 2) Add two new "table" blocks below the existing ones: CONTROL OT
#    (rows 75-80) and AUDITS (rows 82-90). Formatting for each row is
#    cloned from the matching row type further up the sheet, then the
#    text / values are written on top.
# ---------------------------------------------------------------------

# -- row type: section title bar ("Tabla N" / table name) --
$null = $ws.Range("B2:I2").Copy()
$null = $ws.Range("B75").PasteSpecial(-4122)
$null = $ws.Range("B82").PasteSpecial(-4122)

# -- row type: description bar --
$null = $ws.Range("B3:I3").Copy()
$null = $ws.Range("B76").PasteSpecial(-4122)
$null = $ws.Range("B83").PasteSpecial(-4122)

# -- row type: column headers (KEY / COLUMN / TYPE / ...) --
$null = $ws.Range("B4:I4").Copy()
$null = $ws.Range("B77").PasteSpecial(-4122)
$null = $ws.Range("B84").PasteSpecial(-4122)

# -- row type: PK data row (shaded "1" style) --
$null = $ws.Range("B5:I5").Copy()
$null = $ws.Range("B78").PasteSpecial(-4122)
$null = $ws.Range("B85").PasteSpecial(-4122)

# -- row type: FK / alternating data row (shaded "8" style) --
$null = $ws.Range("B6:I6").Copy()
$null = $ws.Range("B79").PasteSpecial(-4122)
$null = $ws.Range("B86").PasteSpecial(-4122)
$null = $ws.Range("B88").PasteSpecial(-4122)
$null = $ws.Range("B90").PasteSpecial(-4122)

# -- row type: plain data row (unshaded "1" style) --
$null = $ws.Range("B7:I7").Copy()
$null = $ws.Range("B80").PasteSpecial(-4122)
$null = $ws.Range("B87").PasteSpecial(-4122)
$null = $ws.Range("B89").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# -- Table "CONTROL OT" --
$ws.Range("B75").Value = "Tabla 9"
$ws.Range("C75").Value = "CONTROL OT"

$ws.Range("B76").Value = "Descripcion"
$ws.Range("C76").Value = "Se crea para tener informacion adicional de cuando se inserta una orden de trabajo"

$ws.Range("B77").Value = "KEY"
$ws.Range("C77").Value = "COLUMN"
$ws.Range("D77").Value = "TYPE"
$ws.Range("E77").Value = "LENGHT"
$ws.Range("F77").Value = "NOT NULL"
$ws.Range("G77").Value = "UNIQUE"
$ws.Range("H77").Value = "DEFAULT"
$ws.Range("I77").Value = "NOTES"

$ws.Range("B78").Value = "PK"
$ws.Range("C78").Value = "CONTROL_ID"
$ws.Range("D78").Value = "INT"
$ws.Range("F78").Value = $true
$ws.Range("G78").Value = $true
$ws.Range("I78").Value = "ID DE LA INSERCION"

$ws.Range("B79").Value = "FK"
$ws.Range("C79").Value = "OR_ID"
$ws.Range("D79").Value = "INT"
$ws.Range("F79").Value = $true
$ws.Range("G79").Value = $true
$ws.Range("I79").Value = "ID DE ORDEN"

$ws.Range("C80").Value = "FECHA_INSERCION"
$ws.Range("D80").Value = "DATETIME"
$ws.Range("F80").Value = $true
$ws.Range("I80").Value = "FECHA EN LA QUE SE INSERTA UNA OT"

# -- Table "AUDITS" --
$ws.Range("B82").Value = "Tabla 10"
$ws.Range("C82").Value = "AUDITS"

$ws.Range("B83").Value = "Descripcion"
$ws.Range("C83").Value = "Se crea para llevar un cotrol de los cambios de las tablas"

$ws.Range("B84").Value = "KEY"
$ws.Range("C84").Value = "COLUMN"
$ws.Range("D84").Value = "TYPE"
$ws.Range("E84").Value = "LENGHT"
$ws.Range("F84").Value = "NOT NULL"
$ws.Range("G84").Value = "UNIQUE"
$ws.Range("H84").Value = "DEFAULT"
$ws.Range("I84").Value = "NOTES"

$ws.Range("B85").Value = "PK"
$ws.Range("C85").Value = "ID_LOG"
$ws.Range("D85").Value = "INT"
$ws.Range("F85").Value = $true
$ws.Range("G85").Value = $true
$ws.Range("I85").Value = "ID DEL LOG"

$ws.Range("C86").Value = "ENTITY"
$ws.Range("D86").Value = "VARCHAR"
$ws.Range("E86").Value = 100
$ws.Range("F86").Value = $true
$ws.Range("I86").Value = "ENTIDAD"

$ws.Range("C87").Value = "ENTITY_ID"
$ws.Range("D87").Value = "INT"
$ws.Range("F87").Value = $true
$ws.Range("I87").Value = "ID DE LA ENTIDAD"

$ws.Range("C88").Value = "OPERATION"
$ws.Range("D88").Value = "VARCHAR"
$ws.Range("E88").Value = 50
$ws.Range("F88").Value = $true
$ws.Range("I88").Value = "OPERACIÓN (INSERT, UPDATE, DELETE)"

$ws.Range("C89").Value = "INSERT_DT"
$ws.Range("D89").Value = "DATETIME"
$ws.Range("F89").Value = $true
$ws.Range("I89").Value = "FECHA DE LA OPERACIÓN"

$ws.Range("C90").Value = "CREATED_BY"
$ws.Range("D90").Value = "VARCHAR"
$ws.Range("E90").Value = 100
$ws.Range("F90").Value = $true
$ws.Range("I90").Value = "QUIEN REALIZO LA OPERACIÓN"

# ---------------------------------------------------------------------
# 3) Merge the title / description header bars of the two new tables,
#    same as every other table block on the sheet.
# ---------------------------------------------------------------------
$null = $ws.Range("C75:I75").Merge()
$null = $ws.Range("C76:I76").Merge()
$null = $ws.Range("C82:I82").Merge()
$null = $ws.Range("C83:I83").Merge()

# ---------------------------------------------------------------------
# 4) Column D widened slightly to fit "DATETIME" / "FECHA_INSERCION".
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 9.6

# ---------------------------------------------------------------------
# 5) Scroll/selection: the sheet no longer freezes on A7 and the final
#    selection lands on K85 (just past the new AUDITS table).
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$null = $ws.Range("K85").Select()
